$wb = $excel.ActiveWorkbook

# ==========================================================================
# Sheet "SparcsN4Login": fill in row 11 (YINV_101) which had blank B:H,
# and insert a new YINV_102 record as row 12 (pushing YINV_103..YINV_107
# down by one row, to rows 13..17).
# ==========================================================================
$ws1 = $wb.Worksheets.Item("SparcsN4Login")

# Make room for the new row at the bottom of the existing block (row 17)
# first, so the row-level formatting already attached to rows 11-16 is left
# completely undisturbed; we then rewrite the row contents directly.
$ws1.Rows.Item(17).Insert()

$loginRows = @(
    @("YINV_101","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_102","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_103","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_104","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_105","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_106","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations"),
    @("YINV_107","admin","Admin","OPR1","CPX11","FCY111","YRD1111","Operations")
)
$r = 11
foreach ($rowVals in $loginRows) {
    $c = 1
    foreach ($v in $rowVals) {
        $ws1.Cells.Item($r, $c).Value = $v
        $c = $c + 1
    }
    $r = $r + 1
}

$ws1.Activate()
$ws1.Range("A11:H17").Select()

# ==========================================================================
# Sheet "N4MobileCommonRoutines": insert a new YINV_102 record as row 4
# (pushing the former rows 4..8 down to 5..9).
# ==========================================================================
$ws10 = $wb.Worksheets.Item("N4MobileCommonRoutines")
$ws10.Rows.Item(4).Insert()

$ws10.Cells.Item(4, 1).Value = "YINV_102"
$ws10.Cells.Item(4, 2).Value = "admin"
$ws10.Cells.Item(4, 3).Value = "Admin"
$ws10.Cells.Item(4, 4).Value = "OPR1"
$ws10.Cells.Item(4, 5).Value = "CPX11"
$ws10.Cells.Item(4, 6).Value = "FCY111"
$ws10.Cells.Item(4, 7).Value = "YRD1111"
$ws10.Cells.Item(4, 8).Value = "Yard Inventory"
$ws10.Cells.Item(4, 9).Value = "Query"
$ws10.Cells.Item(4, 10).Value = "SBSU1234570"

# The inserted row only carries data through column J; clear the blank
# placeholder cells the row-insert leaves behind in K:Q so the row matches
# a genuinely fresh data row with no trailing empty cells.
$ws10.Range("K4:Q4").Clear()

$ws10.Range("A4").Select()

# ==========================================================================
# Sheet "UnitFacilityVisit": append YINV_104..YINV_107 unit rows.
# ==========================================================================
$ws9 = $wb.Worksheets.Item("UnitFacilityVisit")

$ufvRows = @(
    @("YINV_104","SBSU1234570"),
    @("YINV_105","SBSU1234570"),
    @("YINV_106","SBSU1234570"),
    @("YINV_107","SBSU1234570")
)
$r = 3
foreach ($rowVals in $ufvRows) {
    $ws9.Cells.Item($r, 1).Value = $rowVals[0]
    $ws9.Cells.Item($r, 2).Value = $rowVals[1]
    $r = $r + 1
}

$ws9.Range("A1:B6").Select()
